# Apply odds updates to Sheet1 of the workbook per the commit diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Row 5 ---
$ws.Range("G5").Value  = 1.9
$ws.Range("I5").Value  = 3.9
$ws.Range("J5").Value  = 2.63
$ws.Range("S5").Value  = 1.5
$ws.Range("T5").Value  = 2.5
$ws.Range("X5").Value  = 8
$ws.Range("AG5").Value = 9
$ws.Range("AJ5").Value = 41
$ws.Range("AL5").Value = 41
$ws.Range("AN5").Value = 3.75
$ws.Range("AO5").Value = 11
$ws.Range("AT5").Value = 2.5
$ws.Range("AU5").Value = 9
$ws.Range("BA5").Value = 126

# --- Row 17 ---
$ws.Range("G17").Value  = 2
$ws.Range("I17").Value  = 3.5
$ws.Range("J17").Value  = 2.6
$ws.Range("K17").Value  = 2.38
$ws.Range("AH17").Value = 21
$ws.Range("AK17").Value = 26
$ws.Range("AP17").Value = 17
$ws.Range("AY17").Value = 23

# --- Row 24 ---
$ws.Range("J24").Value  = 2.2
$ws.Range("K24").Value  = 2.25
$ws.Range("L24").Value  = 4.75
$ws.Range("S24").Value  = 1.36
$ws.Range("T24").Value  = 2.9
$ws.Range("W24").Value  = 7.5
$ws.Range("Y24").Value  = 8.25
$ws.Range("AK24").Value = 40
$ws.Range("AO24").Value = 8
$ws.Range("AP24").Value = 17
$ws.Range("AQ24").Value = 26
$ws.Range("AR24").Value = 55
$ws.Range("AT24").Value = 2.9
$ws.Range("AW24").Value = 6.3
$ws.Range("AX24").Value = 25
